$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 813
$ws1.Range("F4").Value = 1135
$ws1.Range("F7").Value = 245
$ws1.Range("F14").Value = 12871
$ws1.Range("F16").Value = 5283

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 117

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 813
$ws4.Range("F4").Value = 1135
$ws4.Range("F7").Value = 245
$ws4.Range("F14").Value = 12871
$ws4.Range("F15").Value = 117
$ws4.Range("F18").Value = 5283
